# Corrected excel sheets for application fix issues
#
# Updates the "Summary" and "Repayment schedule" sheets with corrected
# instalment/interest figures, adds a couple of zero-value cells that were
# missing, and refreshes the remembered cell-selection on those two sheets.

$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary")
$wsSchedule  = $wb.Worksheets.Item("Repayment schedule")
$wsEditSched = $wb.Worksheets.Item("Edit Repayment Schedule")

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------

# F2 switches from an integer display to a 2-decimal display (new style),
# and its value is corrected.
$wsSummary.Range("F2").NumberFormat = "#,##0.00"
$wsSummary.Range("F2").Value = 1698.21

# Corrected principal / outstanding figures on row 3.
$wsSummary.Range("A3").Value = 198.47
$wsSummary.Range("E3").Value = 198.47

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------

# New (previously absent) cells - clone the format of their already
# present, identically-styled neighbour so no stray style is created.
$wsSchedule.Range("O2").Copy($wsSchedule.Range("P2"))
$wsSchedule.Range("N3").Copy($wsSchedule.Range("O3"))
$wsSchedule.Range("N4").Copy($wsSchedule.Range("O4"))
$wsSchedule.Range("N5").Copy($wsSchedule.Range("O5"))
$wsSchedule.Range("N6").Copy($wsSchedule.Range("O6"))
$wsSchedule.Range("N7").Copy($wsSchedule.Range("O7"))
$wsSchedule.Range("N8").Copy($wsSchedule.Range("O8"))

$wsSchedule.Range("O3").Value = 0
$wsSchedule.Range("O4").Value = 0
$wsSchedule.Range("O5").Value = 0
$wsSchedule.Range("O6").Value = 0
$wsSchedule.Range("O7").Value = 0
$wsSchedule.Range("O8").Value = 0

# Row 3
$wsSchedule.Range("F3").Value = 848.21
$wsSchedule.Range("G3").NumberFormat = "#,##0.00"
$wsSchedule.Range("G3").Value = 4151.79
$wsSchedule.Range("K3").Value = 900
$wsSchedule.Range("P3").Value = 900

# Row 4
$wsSchedule.Range("F4").Value = 850
$wsSchedule.Range("G4").NumberFormat = "#,##0.00"
$wsSchedule.Range("G4").Value = 3301.79
$wsSchedule.Range("K4").Value = 900
$wsSchedule.Range("P4").Value = 900

# Row 5
$wsSchedule.Range("F5").Value = 850.51
$wsSchedule.Range("G5").NumberFormat = "#,##0.00"
$wsSchedule.Range("G5").Value = 2451.2800000000002
$wsSchedule.Range("H5").Value = 49.49
$wsSchedule.Range("K5").Value = 900
$wsSchedule.Range("P5").Value = 900

# Row 6
$wsSchedule.Range("F6").Value = 875.49
$wsSchedule.Range("G6").NumberFormat = "#,##0.00"
$wsSchedule.Range("G6").Value = 1575.79
$wsSchedule.Range("H6").Value = 24.51
$wsSchedule.Range("K6").Value = 900
$wsSchedule.Range("P6").Value = 900

# Row 7 (G7 keeps its existing style/number format - only the value changes)
$wsSchedule.Range("F7").Value = 884.24
$wsSchedule.Range("G7").Value = 691.55
$wsSchedule.Range("H7").Value = 15.76
$wsSchedule.Range("K7").Value = 900
$wsSchedule.Range("P7").Value = 900

# Row 8
$wsSchedule.Range("F8").Value = 691.55
$wsSchedule.Range("H8").Value = 6.92
$wsSchedule.Range("K8").Value = 698.47
$wsSchedule.Range("P8").Value = 698.47

# ---------------------------------------------------------------------
# Remembered selections
# ---------------------------------------------------------------------

$wsSummary.Range("A7:XFD13").Select()
$wsSchedule.Range("A2:XFD2").Select()

# Restore the workbook's originally active sheet/tab.
$wsEditSched.Activate()
